# The deck ships with two theme parts:
#   theme1.xml ("Integral" / "Red Violet" colours) -> used by the Slide Master
#   theme2.xml ("Office Theme" / "Office" colours)  -> used by the Notes Master
#
# The authored change swaps the two colour palettes between the two themes,
# so the Slide Master (the theme that actually drives the look of every
# slide) now uses the stock "Office" palette instead of "Red Violet".

function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation

# Palette that used to live in theme2.xml ("Office"), keyed by the standard
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @{
    1  = (RGB 0x00 0x00 0x00)   # dk1
    2  = (RGB 0xFF 0xFF 0xFF)   # lt1
    3  = (RGB 0x44 0x54 0x6A)   # dk2
    4  = (RGB 0xE7 0xE6 0xE6)   # lt2
    5  = (RGB 0x5B 0x9B 0xD5)   # accent1
    6  = (RGB 0xED 0x7D 0x31)   # accent2
    7  = (RGB 0xA5 0xA5 0xA5)   # accent3
    8  = (RGB 0xFF 0xC0 0x00)   # accent4
    9  = (RGB 0x44 0x72 0xC4)   # accent5
    10 = (RGB 0x70 0xAD 0x47)   # accent6
    11 = (RGB 0x05 0x63 0xC1)   # hlink
    12 = (RGB 0x95 0x4F 0x72)   # folHlink
}

# Re-colour the Slide Master's theme with the "Office" palette.
$slideTheme = $p.Designs.Item(1).SlideMaster.Theme.ThemeColorScheme
foreach ($idx in 1..12) {
    $slideTheme.Colors($idx).RGB = $officeColors[$idx]
}
